$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '27.084.93'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  -2.58%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.864.97'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  -2.36%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.001'
$ws.Range('D4').Style = 'Normal'
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '305.58'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -2.31%  '
$ws.Range('E6').Value = '  -0.03%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.5156'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  -0.30%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.3759'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -0.57%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.07151'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -1.43%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.8880'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -1.94%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '20.68'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -2.75%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.07545'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -1.51%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '1.866.19'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  -2.49%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '5.304'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -2.84%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '89.49'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -2.74%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '1.002'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -0.07%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.000008474'
$ws.Range('D17').Style = 'Normal'
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '14.08'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  -3.28%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '1.000'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +0.02%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '27.115.11'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -2.64%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '5.008'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -2.90%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '2.090.56'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -3.98%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '10.47'
$ws.Range('D23').Style = 'Normal'
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '6.442'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -3.07%  '
$ws.Range('E25').Value = '  -1.91%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '145.10'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -5.79%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '17.95'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -2.27%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '2.087'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -3.99%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '112.84'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -1.89%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '4.652'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -4.30%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '4.664'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -3.90%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '0.09160'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +0.87%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.05102'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -3.53%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '3.070'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -3.72%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.155'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -6.34%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.7229'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -7.35%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.02039'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -2.82%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '3.088'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +0.20%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '2.485'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -4.72%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '1.073'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -2.02%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.5274'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -5.41%  '
$ws.Range('B42').Value = 'FraxShare'
$ws.Range('C42').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '6.474'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -3.87%  '
$ws.Range('B43').Value = 'Quant'
$ws.Range('C43').Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '116.22'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +0.71%  '
$ws.Range('E44').Value = '  -3.29%  '
$ws.Range('E45').Value = '  -3.71%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '1.0000'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -0.01%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.4613'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -4.34%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '9.959'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -4.80%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '1.562'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -3.57%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '36.54'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -1.28%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '63.40'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -5.39%  '
